# Auto-generated edit script: updates cached Leve profit/price values
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 64703.188
$ws.Range("I28").Value = 85339.336
$ws.Range("J28").Value = 2794.75
$ws.Range("K28").Value = 85339.336
$ws.Range("L28").Value = 2794.75
$ws.Range("M28").Value = -84854.336
$ws.Range("N28").Value = -3764.75

$ws.Range("H64").Value = 5054.56
$ws.Range("J64").Value = 5378.2
$ws.Range("L64").Value = 5378.2
$ws.Range("N64").Value = -5874.2

$ws.Range("H67").Value = 5054.56
$ws.Range("J67").Value = 5378.2
$ws.Range("L67").Value = 5378.2
$ws.Range("N67").Value = -7094.2

$ws.Range("I74").Value = 4250
$ws.Range("J74").Value = 6333.3335
$ws.Range("K74").Value = 4250
$ws.Range("L74").Value = 6333.3335
$ws.Range("M74").Value = -3314
$ws.Range("N74").Value = -8205.333500000001

$ws.Range("I77").Value = 4250
$ws.Range("J77").Value = 6333.3335
$ws.Range("K77").Value = 21250
$ws.Range("L77").Value = 31666.6675
$ws.Range("M77").Value = -16570
$ws.Range("N77").Value = -41026.6675

$ws.Range("H125").Value = 947.25
$ws.Range("J125").Value = 947.25
$ws.Range("L125").Value = 8525.25
$ws.Range("N125").Value = -13445.25

$ws.Range("H127").Value = 2113.1428
$ws.Range("J127").Value = 2519.6
$ws.Range("L127").Value = 7558.799999999999
$ws.Range("N127").Value = -17478.8

$ws.Range("H135").Value = 640.54346
$ws.Range("I135").Value = 417.89474
$ws.Range("K135").Value = 3761.05266
$ws.Range("M135").Value = -1226.05266

$ws.Range("H137").Value = 5505.923
$ws.Range("I137").Value = 4740.273
$ws.Range("J137").Value = 6067.4
$ws.Range("K137").Value = 14220.819
$ws.Range("L137").Value = 18202.2
$ws.Range("M137").Value = -11670.819
$ws.Range("N137").Value = -23302.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16057.023
$ws.Range("I32").Value = 8324.831
$ws.Range("J32").Value = 50368.625
$ws.Range("K32").Value = 8324.831
$ws.Range("L32").Value = 50368.625
$ws.Range("M32").Value = -8037.831
$ws.Range("N32").Value = -50942.625

$ws.Range("H45").Value = 328391.34
$ws.Range("I45").Value = 428957.22
$ws.Range("J45").Value = 1552.25
$ws.Range("K45").Value = 428957.22
$ws.Range("L45").Value = 1552.25
$ws.Range("M45").Value = -428580.22
$ws.Range("N45").Value = -2306.25

$ws.Range("H63").Value = 54549630
$ws.Range("I63").Value = 83335784
$ws.Range("J63").Value = 20006242
$ws.Range("K63").Value = 83335784
$ws.Range("L63").Value = 20006242
$ws.Range("M63").Value = -83335098
$ws.Range("N63").Value = -20007614

$ws.Range("H66").Value = 54549630
$ws.Range("I66").Value = 83335784
$ws.Range("J66").Value = 20006242
$ws.Range("K66").Value = 416678920
$ws.Range("L66").Value = 100031210
$ws.Range("M66").Value = -416675488
$ws.Range("N66").Value = -100038074

$ws.Range("H88").Value = 23812488
$ws.Range("J88").Value = 3458.1667
$ws.Range("L88").Value = 3458.1667
$ws.Range("N88").Value = -4270.1667

$ws.Range("H91").Value = 23812488
$ws.Range("J91").Value = 3458.1667
$ws.Range("L91").Value = 3458.1667
$ws.Range("N91").Value = -6266.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6527.7144
$ws.Range("I86").Value = 2850
$ws.Range("J86").Value = 7998.8
$ws.Range("K86").Value = 2850
$ws.Range("L86").Value = 7998.8
$ws.Range("M86").Value = -1727
$ws.Range("N86").Value = -10244.8

$ws.Range("H89").Value = 6527.7144
$ws.Range("I89").Value = 2850
$ws.Range("J89").Value = 7998.8
$ws.Range("K89").Value = 14250
$ws.Range("L89").Value = 39994
$ws.Range("M89").Value = -8634
$ws.Range("N89").Value = -51226

$ws.Range("H134").Value = 2127.6099
$ws.Range("I134").Value = 1993.3
$ws.Range("K134").Value = 5979.9
$ws.Range("M134").Value = -3444.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3619.8262
$ws.Range("I31").Value = 2063.3157
$ws.Range("J31").Value = 4715.148
$ws.Range("K31").Value = 2063.3157
$ws.Range("L31").Value = 4715.148
$ws.Range("M31").Value = -1768.3157
$ws.Range("N31").Value = -5305.148

$ws.Range("H34").Value = 3619.8262
$ws.Range("I34").Value = 2063.3157
$ws.Range("J34").Value = 4715.148
$ws.Range("K34").Value = 2063.3157
$ws.Range("L34").Value = 4715.148
$ws.Range("M34").Value = -1861.3157
$ws.Range("N34").Value = -5119.148

$ws.Range("H58").Value = 360720.03
$ws.Range("I58").Value = 1856.8572
$ws.Range("K58").Value = 1856.8572
$ws.Range("M58").Value = -1653.8572

$ws.Range("H94").Value = 1246.5714
$ws.Range("J94").Value = 1019.3333
$ws.Range("L94").Value = 1019.3333
$ws.Range("N94").Value = -1921.3333

$ws.Range("H132").Value = 403111.62
$ws.Range("I132").Value = 3350
$ws.Range("K132").Value = 10050
$ws.Range("M132").Value = -7520

$ws.Range("H136").Value = 360720.03
$ws.Range("I136").Value = 1856.8572
$ws.Range("K136").Value = 5570.571599999999
$ws.Range("M136").Value = -3020.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36517890
$ws.Range("I4").Value = 37707692
$ws.Range("K4").Value = 113123076
$ws.Range("M4").Value = -113122964

$ws.Range("H57").Value = 489.5
$ws.Range("I57").Value = 489.5
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 1468.5
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -909.5
$ws.Range("N57").ClearContents()

$ws.Range("H109").Value = 3111.8572
$ws.Range("I109").Value = 2062.889
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 6188.667
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -5148.667
$ws.Range("N109").Value = -17080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5023.243
$ws.Range("I132").Value = 4474.76
$ws.Range("J132").Value = 6165.9165
$ws.Range("K132").Value = 13424.28
$ws.Range("L132").Value = 18497.7495
$ws.Range("M132").Value = -10894.28
$ws.Range("N132").Value = -23557.7495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 71433240
$ws.Range("I7").Value = 125002920
$ws.Range("K7").Value = 125002920
$ws.Range("M7").Value = -125002808

$ws.Range("H16").Value = 734.55554
$ws.Range("I16").Value = 758.13336
$ws.Range("J16").Value = 616.6667
$ws.Range("K16").Value = 758.13336
$ws.Range("L16").Value = 616.6667
$ws.Range("M16").Value = -588.13336
$ws.Range("N16").Value = -956.6667

$ws.Range("H22").Value = 613.4
$ws.Range("I22").Value = 621.2857
$ws.Range("J22").Value = 595
$ws.Range("K22").Value = 621.2857
$ws.Range("L22").Value = 595
$ws.Range("M22").Value = -326.2857
$ws.Range("N22").Value = -1185

$ws.Range("H27").Value = 613.4
$ws.Range("I27").Value = 621.2857
$ws.Range("J27").Value = 595
$ws.Range("K27").Value = 621.2857
$ws.Range("L27").Value = 595
$ws.Range("M27").Value = -514.2857
$ws.Range("N27").Value = -809

$ws.Range("H53").Value = 7951
$ws.Range("J53").Value = 7951
$ws.Range("L53").Value = 7951
$ws.Range("N53").Value = -8987

$ws.Range("H55").Value = 1317.6
$ws.Range("I55").Value = 638.4
$ws.Range("J55").Value = 1996.8
$ws.Range("K55").Value = 638.4
$ws.Range("L55").Value = 1996.8
$ws.Range("M55").Value = -465.4
$ws.Range("N55").Value = -2342.8

$ws.Range("H68").Value = 4441.0835
$ws.Range("J68").Value = 8740
$ws.Range("L68").Value = 8740
$ws.Range("N68").Value = -10238

$ws.Range("H71").Value = 4441.0835
$ws.Range("J71").Value = 8740
$ws.Range("L71").Value = 43700
$ws.Range("N71").Value = -51188

$ws.Range("H94").Value = 40001
$ws.Range("J94").Value = 40001
$ws.Range("L94").Value = 40001
$ws.Range("N94").Value = -41353

$ws.Range("H126").Value = 71433240
$ws.Range("I126").Value = 125002920
$ws.Range("K126").Value = 375008760
$ws.Range("M126").Value = -375006290

$ws.Range("H132").Value = 154234.66
$ws.Range("I132").Value = 273708.34
$ws.Range("J132").Value = 6883.7666
$ws.Range("K132").Value = 821125.02
$ws.Range("L132").Value = 20651.2998
$ws.Range("M132").Value = -818595.02
$ws.Range("N132").Value = -25711.2998

$ws.Range("H136").Value = 6909.3125
$ws.Range("I136").Value = 6996.375
$ws.Range("K136").Value = 20989.125
$ws.Range("M136").Value = -18439.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9611.468999999999
$ws.Range("I136").Value = 11924.353
$ws.Range("K136").Value = 35773.05899999999
$ws.Range("M136").Value = -33223.05899999999
